$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel's "looks like
# a number" auto-coercion (which would turn "002292" into 2292, or "0.39"
# into a float). We stash a text-formula result in a scratch cell, copy it,
# and paste-special *values only* into the destination - the pasted value
# keeps its string type without leaving a formula behind or minting a
# quotePrefix style.
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $scratch = $range.Worksheet.Range("ZZ1000")
    $scratch.Formula = "=" + '"' + $text + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# 1. Duplicate the existing "2021-Q3" sheet, placing the copy right after it.
#    The copy will keep the old 2021-Q3 figures untouched and is renamed
#    back to "2021-Q3" below; the original sheet object is recycled to hold
#    the new 2022-Q3 figures (so it keeps rId2 / sheetId 2).
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsOld2021 = $wb.Worksheets.Item("2021-Q3")
$wsOld2021.Copy($null, $wsOld2021)
$wsDup2021 = $wb.Worksheets.Item("2021-Q3 (2)")

# ---------------------------------------------------------------------------
# 2. Turn the original sheet into the "2022-Q3" data sheet.
# ---------------------------------------------------------------------------
$ws2022 = $wsOld2021
$ws2022.Cells.Clear()

Set-TextValue $ws2022.Range("B1") "基金代码"
Set-TextValue $ws2022.Range("C1") "基金名称"
Set-TextValue $ws2022.Range("D1") "基金规模"
Set-TextValue $ws2022.Range("E1") "股票总仓位"
Set-TextValue $ws2022.Range("F1") "仓位占比"
Set-TextValue $ws2022.Range("G1") "持有市值(亿元)"
Set-TextValue $ws2022.Range("H1") "仓位排名"

$ws2022.Range("A2").Value = 0
Set-TextValue $ws2022.Range("B2") "002292"
Set-TextValue $ws2022.Range("C2") "诺安益鑫灵活配置混合A"
Set-TextValue $ws2022.Range("D2") "0.39"
Set-TextValue $ws2022.Range("E2") "61.16"
Set-TextValue $ws2022.Range("F2") "3.44"
Set-TextValue $ws2022.Range("G2") "0.0134"
$ws2022.Range("H2").Value = 7

$ws2022.Range("A3").Value = 1
Set-TextValue $ws2022.Range("B3") "001932"
Set-TextValue $ws2022.Range("C3") "国寿安保灵活优选混合"
Set-TextValue $ws2022.Range("D3") "0.13"
Set-TextValue $ws2022.Range("E3") "34.76"
Set-TextValue $ws2022.Range("F3") "1.11"
Set-TextValue $ws2022.Range("G3") "0.0014"
$ws2022.Range("H3").Value = 8

$ws2022.Range("A4").Value = 2
Set-TextValue $ws2022.Range("B4") "014550"
Set-TextValue $ws2022.Range("C4") "诺安益鑫灵活配置混合C"
Set-TextValue $ws2022.Range("D4") "0.02"
Set-TextValue $ws2022.Range("E4") "61.16"
Set-TextValue $ws2022.Range("F4") "3.44"
Set-TextValue $ws2022.Range("G4") "0.0007"
$ws2022.Range("H4").Value = 7

# Match the bold/bordered header style already used on the "总计" sheet
# (reuse its style instead of minting a new font/border combo).
$wsTotal.Range("B1").Copy()
$ws2022.Range("B1:H1").PasteSpecial(-4122)
$ws2022.Range("A2:A4").PasteSpecial(-4122)

$ws2022.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 3. Rename the duplicated sheet back to "2021-Q3" (keeps its original data).
# ---------------------------------------------------------------------------
$wsDup2021.Name = "2021-Q3"

# ---------------------------------------------------------------------------
# 4. Update the "总计" summary sheet: the former row 2 ("2021-Q3" totals)
#    moves down to row 3, and row 2 now holds the new "2022-Q3" totals.
# ---------------------------------------------------------------------------
$wsTotal.Range("A3").Value = 1
Set-TextValue $wsTotal.Range("B3") "2021-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.13

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

Set-TextValue $wsTotal.Range("B2") "2022-Q3"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.02
